$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns before column B (old B:E shift right to E:H)
$ws.Columns("B:D").Insert()

# New header values for the inserted/shifted columns in row 1
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the three newly inserted data columns (B, C, D) for every existing data row
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B$r").Value = "UN"
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
}

# Insert 2 new rows at the bottom for the new analyst group
$ws.Rows("28:29").Insert()

$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
